$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dept Data")

# Add two new columns capturing bank/account-type info for each customer:
# "banktype" (F) and "AccountType" (G)
$ws.Range("F1").Value = "banktype"
$ws.Range("G1").Value = "AccountType"

$data = @(
    @("Inter-bank", "Savings"),
    @("Inter-bank", "Savings"),
    @("Intra-bank", "Savings"),
    @("Inter-bank", "Current"),
    @("Inter-bank", "Savings"),
    @("Intra-bank", "Current")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $data[$i][0]
    $ws.Cells.Item($row, 7).Value = $data[$i][1]
}
